$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.008627148837446561
$ws.Range("J2").Value = 0.008627148837446563
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2317656666666667
$ws.Range("N2").Value = 0.6952970000000001
$ws.Range("O2").Value = 0.04497992035421015
$ws.Range("P2").Value = 0.04579239628237638
$ws.Range("Q2").Value = 0.004024224525555556
$ws.Range("R2").Value = 0.03621802073
$ws.Range("S2").Value = 0.000388048467592263
$ws.Range("T2").Value = 0.0003950578183513957
$ws.Range("I3").Value = 0.008627148837446561
$ws.Range("J3").Value = 0.008627148837446563
$ws.Range("O3").Value = 0.1919731617572378
$ws.Range("P3").Value = 0.1954407884571841
$ws.Range("S3").Value = 0.001656181039274895
$ws.Range("T3").Value = 0.001686096770928035
$ws.Range("I4").Value = 0.008627148837446561
$ws.Range("J4").Value = 0.008627148837446563
$ws.Range("M4").Value = 2.312753666666667
$ws.Range("N4").Value = 6.938261000000001
$ws.Range("O4").Value = 0.4488476538468058
$ws.Range("P4").Value = 0.4569552252095968
$ws.Range("Q4").Value = 0.04015711283222222
$ws.Range("R4").Value = 0.36141401549
$ws.Range("S4").Value = 0.003872275515075087
$ws.Range("T4").Value = 0.003942220739932106
$ws.Range("I5").Value = 0.008627148837446561
$ws.Range("J5").Value = 0.008627148837446563
$ws.Range("M5").Value = 0.274264
$ws.Range("N5").Value = 0.548528
$ws.Range("O5").Value = 0.05322778413840601
$ws.Range("P5").Value = 0.03612616126343038
$ws.Range("Q5").Value = 0.004762137253333333
$ws.Range("R5").Value = 0.02857282352
$ws.Range("S5").Value = 0.0004592040160495059
$ws.Range("T5").Value = 0.0003116657701452104
$ws.Range("I6").Value = 0.008627148837446561
$ws.Range("J6").Value = 0.008627148837446563
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.344694
$ws.Range("N6").Value = 4.034082
$ws.Range("O6").Value = 0.2609714799033402
$ws.Range("P6").Value = 0.2656854287874124
$ws.Range("Q6").Value = 0.02334837015333333
$ws.Range("R6").Value = 0.21013533138
$ws.Range("S6").Value = 0.002251439799454811
$ws.Range("T6").Value = 0.002292107738089816
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.023931333333333
$ws.Range("H7").Value = 3.071794
$ws.Range("I7").Value = 0.5087507014009469
$ws.Range("J7").Value = 0.5087507014009469
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2317656666666667
$ws.Range("N7").Value = 0.6952970000000001
$ws.Range("O7").Value = 0.04497992035421015
$ws.Range("P7").Value = 0.04579239628237638
$ws.Range("Q7").Value = 0.2373121280908889
$ws.Range("R7").Value = 2.135809152818
$ws.Range("S7").Value = 0.02288356602916315
$ws.Range("T7").Value = 0.0232969137274891
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.023931333333333
$ws.Range("H8").Value = 3.071794
$ws.Range("I8").Value = 0.5087507014009469
$ws.Range("J8").Value = 0.5087507014009469
$ws.Range("O8").Value = 0.1919731617572378
$ws.Range("P8").Value = 0.1954407884571841
$ws.Range("Q8").Value = 1.012842156993333
$ws.Range("R8").Value = 9.115579412939999
$ws.Range("S8").Value = 0.09766648069415217
$ws.Range("T8").Value = 0.09943063820994649
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.023931333333333
$ws.Range("H9").Value = 3.071794
$ws.Range("I9").Value = 0.5087507014009469
$ws.Range("J9").Value = 0.5087507014009469
$ws.Range("M9").Value = 2.312753666666667
$ws.Range("N9").Value = 6.938261000000001
$ws.Range("O9").Value = 0.4488476538468058
$ws.Range("P9").Value = 0.4569552252095968
$ws.Range("Q9").Value = 2.368100945581556
$ws.Range("R9").Value = 21.312908510234
$ws.Range("S9").Value = 0.2283515587167319
$ws.Range("T9").Value = 0.23247629133421
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.023931333333333
$ws.Range("H10").Value = 3.071794
$ws.Range("I10").Value = 0.5087507014009469
$ws.Range("J10").Value = 0.5087507014009469
$ws.Range("M10").Value = 0.274264
$ws.Range("N10").Value = 0.548528
$ws.Range("O10").Value = 0.05322778413840601
$ws.Range("P10").Value = 0.03612616126343038
$ws.Range("Q10").Value = 0.2808275032053333
$ws.Range("R10").Value = 1.684965019232
$ws.Range("S10").Value = 0.02707967251443225
$ws.Range("T10").Value = 0.01837920988169392
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.023931333333333
$ws.Range("H11").Value = 3.071794
$ws.Range("I11").Value = 0.5087507014009469
$ws.Range("J11").Value = 0.5087507014009469
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.344694
$ws.Range("N11").Value = 4.034082
$ws.Range("O11").Value = 0.2609714799033402
$ws.Range("P11").Value = 0.2656854287874124
$ws.Range("Q11").Value = 1.376874320345333
$ws.Range("R11").Value = 12.391868883108
$ws.Range("S11").Value = 0.1327694234464675
$ws.Range("T11").Value = 0.1351676482476074
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9713440000000001
$ws.Range("H12").Value = 2.914032
$ws.Range("I12").Value = 0.4826221497616066
$ws.Range("J12").Value = 0.4826221497616065
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2317656666666667
$ws.Range("N12").Value = 0.6952970000000001
$ws.Range("O12").Value = 0.04497992035421015
$ws.Range("P12").Value = 0.04579239628237638
$ws.Range("Q12").Value = 0.2251241897226667
$ws.Range("R12").Value = 2.026117707504
$ws.Range("S12").Value = 0.02170830585745475
$ws.Range("T12").Value = 0.02210042473653588
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9713440000000001
$ws.Range("H13").Value = 2.914032
$ws.Range("I13").Value = 0.4826221497616066
$ws.Range("J13").Value = 0.4826221497616065
$ws.Range("O13").Value = 0.1919731617572378
$ws.Range("P13").Value = 0.1954407884571841
$ws.Range("Q13").Value = 0.9608243444800001
$ws.Range("R13").Value = 8.64741910032
$ws.Range("S13").Value = 0.09265050002381074
$ws.Range("T13").Value = 0.09432405347630955
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.9713440000000001
$ws.Range("H14").Value = 2.914032
$ws.Range("I14").Value = 0.4826221497616066
$ws.Range("J14").Value = 0.4826221497616065
$ws.Range("M14").Value = 2.312753666666667
$ws.Range("N14").Value = 6.938261000000001
$ws.Range("O14").Value = 0.4488476538468058
$ws.Range("P14").Value = 0.4569552252095968
$ws.Range("Q14").Value = 2.246479397594667
$ws.Range("R14").Value = 20.218314578352
$ws.Range("S14").Value = 0.2166238196149989
$ws.Range("T14").Value = 0.2205367131354546
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.9713440000000001
$ws.Range("H15").Value = 2.914032
$ws.Range("I15").Value = 0.4826221497616066
$ws.Range("J15").Value = 0.4826221497616065
$ws.Range("M15").Value = 0.274264
$ws.Range("N15").Value = 0.548528
$ws.Range("O15").Value = 0.05322778413840601
$ws.Range("P15").Value = 0.03612616126343038
$ws.Range("Q15").Value = 0.266404690816
$ws.Range("R15").Value = 1.598428144896
$ws.Range("S15").Value = 0.02568890760792425
$ws.Range("T15").Value = 0.01743528561159124
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.9713440000000001
$ws.Range("H16").Value = 2.914032
$ws.Range("I16").Value = 0.4826221497616066
$ws.Range("J16").Value = 0.4826221497616065
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.344694
$ws.Range("N16").Value = 4.034082
$ws.Range("O16").Value = 0.2609714799033402
$ws.Range("P16").Value = 0.2656854287874124
$ws.Range("Q16").Value = 1.306160448736
$ws.Range("R16").Value = 11.755444038624
$ws.Range("S16").Value = 0.125950616657418
$ws.Range("T16").Value = 0.1282256728017152
